$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 13: "2021年" summary row, appended after existing 2010-2020 data (rows 2-12).
$ws.Range("A13").Value = "2021年"

# Match the formatting used by the other year cells in column A (bold, centered,
# bordered style) by copying the format from the cell directly above (A12).
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$values = @(101, 1211, 126, 360, 888, 270, 2025, 762, 2, 3218, 2659, 445, 153, 12220)
for ($i = 0; $i -lt $values.Length; $i++) {
    $col = 2 + $i
    $ws.Cells.Item(13, $col).Value = $values[$i]
}
